$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "e010 Time Check" description in cell B11 with the revised
# wording (added the "Roll 1D/2" sentence and reworked the closing
# paragraph about the Time Table / ammo expended).
$lines = @(
    '<Bold>e010 Time Check</Bold> ',
    '<InlineUIContainer><Button Content=''r4.3'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>  ',
    '<LineBreak/><LineBreak/>',
    'Determine sunrise and sunset for current month using the <InlineUIContainer><Button Content=''Time'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer> Table. Roll 1D/2  on the ',
    '<InlineUIContainer><Button Content=''Time'' FontFamily=''Courier New''  FontSize=''12''></Button></InlineUIContainer>  Table. The Time Table also provides the timed used for each action take. Additionally, the same die roll is used to determine the ammo expended:  ',
    '<InlineUIContainer><Image Name=''DieRoll'' Height=''21'' Width=''21'' > </Image></InlineUIContainer>',
    '<LineBreak/><LineBreak/>'
)

$newText = [string]::Join("`r`n", $lines)
$ws.Range("B11").Value = $newText

# The extra wrapped line pushes row 11 taller (one default row height
# taller: 99.85 -> 114.15).
$ws.Range("A11").EntireRow.RowHeight = 114.15

# Reflect where the user left the selection after editing the cell.
$ws.Range("B9").Select()
